$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.438.63"
$ws.Range("E2").Value = "  -3.80%  "

$ws.Range("D3").Value = "1.769.75"
$ws.Range("E3").Value = "  -3.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.58"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4295"
$ws.Range("E7").Value = "  +0.84%  "

$ws.Range("E8").Value = "  +1.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07223"
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8490"
$ws.Range("E10").Value = "  -1.73%  "

$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").Value = "1.791.00"
$ws.Range("E12").Value = "  -8.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.428"
$ws.Range("E13").Value = "  -0.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.233"
$ws.Range("E14").Value = "  -2.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06899"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.006"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.33"
$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008686"
$ws.Range("E18").Value = "  -2.85%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.04"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("D21").Value = "26.434.17"
$ws.Range("E21").Value = "  -5.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.101"
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.28"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("D24").Value = "1.976.84"
$ws.Range("E24").Value = "  -4.62%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.879"
$ws.Range("E25").Value = "  -5.58%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.21"
$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.11"
$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.093"
$ws.Range("E28").Value = "  -1.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.53"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.728"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08966"
$ws.Range("E31").Value = "  +0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7223"
$ws.Range("E32").Value = "  -3.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.112"
$ws.Range("E33").Value = "  -0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.320"
$ws.Range("E34").Value = "  -4.94%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.744"
$ws.Range("E36").Value = "  -7.77%  "

$ws.Range("E37").Value = "  -0.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05162"
$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01892"
$ws.Range("E39").Value = "  -1.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4923"
$ws.Range("E40").Value = "  -3.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1601"
$ws.Range("E41").Value = "  -3.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.575"
$ws.Range("E42").Value = "  -7.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.236"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.001"
$ws.Range("E44").Value = "  -4.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.78"
$ws.Range("E45").Value = "  -1.61%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.11"
$ws.Range("E47").Value = "  -3.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06199"
$ws.Range("E48").Value = "  -4.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4493"
$ws.Range("E49").Value = "  -3.90%  "

$ws.Range("E50").Value = "  -1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.738"
$ws.Range("E51").Value = "  +2.35%  "
